# Multiple tables using "getJoins()". Still need to do "WHERE xxx.yyy = www.zzz" (not working yet)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Add two new data rows to Sheet1 (row 3 and row 4)
$ws1.Range("A3").Value = 5
$ws1.Range("B3").Value = $false
$ws1.Range("C3").Value = """And Dad"""
$ws1.Range("D3").Value = 3

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = $true
$ws1.Range("C4").Value = """foobar"""
$ws1.Range("D4").Value = 6

# Sheet1 becomes the active/selected sheet (was Sheet2 before), with E11 selected
$ws1.Activate()
$ws1.Range("E11").Select()
